$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing sensor readings (columns C:H, rows 2-20) down by one row
# into rows 3-21, making room for a newly captured reading at row 2.
$src = $ws.Range("C2:H20")
$dst = $ws.Range("C3:H21")
$dst.Value2 = $src.Value2

# Write the new reading into row 2
$ws.Cells.Item(2,3).Value2 = -3.540287351608276
$ws.Cells.Item(2,4).Value2 = 5.382533311843872
$ws.Cells.Item(2,5).Value2 = -1.780441856384277
$ws.Cells.Item(2,6).Value2 = -0.0032070425804704
$ws.Cells.Item(2,7).Value2 = 0.0187841057777404
$ws.Cells.Item(2,8).Value2 = 0.0213802829384803

# Drop the oldest reading, which is now pushed beyond the tracked window
$ws.Rows.Item(22).Delete()
